# Refresh the cryptos price/volume snapshot (GitHub Actions scheduled update).
# Price/percent columns are stored as literal text (e.g. "591.60", "  +3.84%  "),
# so numeric-looking prices are written through a helper that forces the Text
# number format for the duration of the write, then restores whatever format
# the cell had before -- this stops Excel from silently re-typing "591.60" as
# the number 591.6 (dropping the trailing zero) while leaving no lasting format change.

function Set-TextValue($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = '@'
    $cell.Value = $text
    $cell.NumberFormat = $originalFormat
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Cells.Item(2, 4).Value = '66.647.94'
$ws.Cells.Item(2, 5).Value = '  +4.51%  '

# Row 3: Ethereum
$ws.Cells.Item(3, 4).Value = '3.486.64'
$ws.Cells.Item(3, 5).Value = '  +2.78%  '

# Row 4: TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.03%  '

# Row 5: BNB
Set-TextValue $ws 5 4 '591.60'
$ws.Cells.Item(5, 5).Value = '  +3.84%  '

# Row 6: Solana
Set-TextValue $ws 6 4 '169.31'
$ws.Cells.Item(6, 5).Value = '  +5.08%  '

# Row 7: USDC
Set-TextValue $ws 7 4 '1.00'
$ws.Cells.Item(7, 5).Value = '  -0.03%  '

# Row 8: LidoStakedEther
$ws.Cells.Item(8, 4).Value = '3.485.03'
$ws.Cells.Item(8, 5).Value = '  +2.76%  '

# Row 9: XRP
Set-TextValue $ws 9 4 '0.593'
$ws.Cells.Item(9, 5).Value = '  +8.84%  '

# Row 10: Toncoin
Set-TextValue $ws 10 4 '7.33'

# Row 11: Dogecoin
Set-TextValue $ws 11 4 '0.127'
$ws.Cells.Item(11, 5).Value = '  +7.45%  '

# Row 12: Cardano
Set-TextValue $ws 12 4 '0.437'
$ws.Cells.Item(12, 5).Value = '  +4.44%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Cells.Item(13, 4).Value = '4.092.52'
$ws.Cells.Item(13, 5).Value = '  +2.89%  '

# Row 14: TRON
$ws.Cells.Item(14, 5).Value = '  -0.19%  '

# Row 15: Avalanche
Set-TextValue $ws 15 4 '28.14'
$ws.Cells.Item(15, 5).Value = '  +4.92%  '

# Row 16: ShibaInu
Set-TextValue $ws 16 4 '0.0000178'
$ws.Cells.Item(16, 5).Value = '  +4.15%  '

# Row 17: WrappedBTC
$ws.Cells.Item(17, 4).Value = '66.673.58'
$ws.Cells.Item(17, 5).Value = '  +4.52%  '

# Row 18: WrappedEther
$ws.Cells.Item(18, 4).Value = '3.502.03'
$ws.Cells.Item(18, 5).Value = '  +3.14%  '

# Row 19: Polkadot
Set-TextValue $ws 19 4 '6.31'
$ws.Cells.Item(19, 5).Value = '  +3.70%  '

# Row 20: Chainlink
Set-TextValue $ws 20 4 '14.04'
$ws.Cells.Item(20, 5).Value = '  +4.18%  '

# Row 21: BitcoinCash
Set-TextValue $ws 21 4 '392.33'
$ws.Cells.Item(21, 5).Value = '  +4.42%  '

# Row 22: Uniswap
Set-TextValue $ws 22 4 '7.91'
$ws.Cells.Item(22, 5).Value = '  +2.06%  '

# Row 23: Litecoin
Set-TextValue $ws 23 4 '73.09'
$ws.Cells.Item(23, 5).Value = '  +4.53%  '

# Row 24: Dai
$ws.Cells.Item(24, 5).Value = '  -0.17%  '

# Row 25: Polygon
Set-TextValue $ws 25 4 '0.534'
$ws.Cells.Item(25, 5).Value = '  +4.76%  '

# Row 26: PEPE
$ws.Cells.Item(26, 5).Value = '  +7.18%  '

# Row 27: InternetComputer(DFINITY)
Set-TextValue $ws 27 4 '10.24'
$ws.Cells.Item(27, 5).Value = '  +7.67%  '

# Row 28: Kaspa
Set-TextValue $ws 28 4 '0.181'
$ws.Cells.Item(28, 5).Value = '  +1.60%  '

# Row 29: Binance-PegBSC-USD
$ws.Cells.Item(29, 5).Value = '  -0.19%  '

# Row 30: NEARProtocol
Set-TextValue $ws 30 4 '6.33'
$ws.Cells.Item(30, 5).Value = '  +4.87%  '

# Row 31: Fetch.AI
$ws.Cells.Item(31, 5).Value = '  +5.46%  '

# Row 32: PancakeSwap
$ws.Cells.Item(32, 5).Value = '  +3.83%  '

# Row 33: EthereumClassic
Set-TextValue $ws 33 4 '23.52'
$ws.Cells.Item(33, 5).Value = '  +3.61%  '

# Row 34: Aptos
$ws.Cells.Item(34, 5).Value = '  +5.92%  '

# Row 35: USDe
$ws.Cells.Item(35, 5).Value = '  +0.12%  '

# Row 36: ImmutableX
$ws.Cells.Item(36, 5).Value = '  +10.34%  '

# Row 37: Monero
Set-TextValue $ws 37 4 '161.55'
$ws.Cells.Item(37, 5).Value = '  +1.38%  '

# Row 38: Mantle
Set-TextValue $ws 38 4 '0.900'
$ws.Cells.Item(38, 5).Value = '  +5.05%  '

# Row 39: Stacks
$ws.Cells.Item(39, 5).Value = '  +7.41%  '

# Row 40: RenderToken
Set-TextValue $ws 40 4 '6.76'
$ws.Cells.Item(40, 5).Value = '  +5.62%  '

# Row 41: Hedera
$ws.Cells.Item(41, 5).Value = '  +3.77%  '

# Row 42: EnergySwap (was Filecoin)
$ws.Cells.Item(42, 2).Value = 'EnergySwap'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws 42 4 '26.60'
$ws.Cells.Item(42, 5).Value = '  +3.71%  '

# Row 43: Filecoin (was EnergySwap)
$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws 43 4 '4.65'
$ws.Cells.Item(43, 5).Value = '  +7.15%  '

# Row 44: InjectiveProtocol
Set-TextValue $ws 44 4 '26.78'
$ws.Cells.Item(44, 5).Value = '  +3.07%  '

# Row 45: OKB
Set-TextValue $ws 45 4 '43.20'
$ws.Cells.Item(45, 5).Value = '  +1.37%  '

# Row 46: Maker
$ws.Cells.Item(46, 4).Value = '2.762.71'
$ws.Cells.Item(46, 5).Value = '  +1.05%  '

# Row 47: VeChain
Set-TextValue $ws 47 4 '0.0312'
$ws.Cells.Item(47, 5).Value = '  +2.64%  '

# Row 48: dogwifhat
Set-TextValue $ws 48 4 '2.48'
$ws.Cells.Item(48, 5).Value = '  +3.75%  '

# Row 49: Bittensor
Set-TextValue $ws 49 4 '346.36'
$ws.Cells.Item(49, 5).Value = '  +5.86%  '

# Row 50: ONDO
$ws.Cells.Item(50, 5).Value = '  +5.66%  '

# Row 51: Arweave (was SuiNetwork)
$ws.Cells.Item(51, 2).Value = 'Arweave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue $ws 51 4 '34.00'
$ws.Cells.Item(51, 5).Value = '  +14.16%  '
